$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 96
$ws.Cells.Item(5, 10).Value = 82.42856999999999
$ws.Cells.Item(5, 12).Value = 82.42856999999999
$ws.Cells.Item(5, 14).Value = -312.42857

$ws.Cells.Item(12, 8).Value = 115.77778
$ws.Cells.Item(12, 9).Value = 99.333336
$ws.Cells.Item(12, 10).Value = 198
$ws.Cells.Item(12, 11).Value = 99.333336
$ws.Cells.Item(12, 12).Value = 198
$ws.Cells.Item(12, 13).Value = 70.666664
$ws.Cells.Item(12, 14).Value = -538

$ws.Cells.Item(17, 8).Value = 1176.3914
$ws.Cells.Item(17, 9).Value = 499
$ws.Cells.Item(17, 10).Value = 1278
$ws.Cells.Item(17, 11).Value = 1497
$ws.Cells.Item(17, 12).Value = 3834
$ws.Cells.Item(17, 13).Value = -1329
$ws.Cells.Item(17, 14).Value = -4170

$ws.Cells.Item(96, 8).Value = 349.26666
$ws.Cells.Item(96, 9).Value = 288.3
$ws.Cells.Item(96, 10).Value = 471.2
$ws.Cells.Item(96, 11).Value = 864.9000000000001
$ws.Cells.Item(96, 12).Value = 1413.6
$ws.Cells.Item(96, 13).Value = 508.0999999999999
$ws.Cells.Item(96, 14).Value = -4159.6

$ws.Cells.Item(100, 8).Value = 33334974
$ws.Cells.Item(100, 9).Value = 1600.4546
$ws.Cells.Item(100, 11).Value = 1600.4546
$ws.Cells.Item(100, 13).Value = -1059.4546

$ws.Cells.Item(121, 8).Value = 987.2
$ws.Cells.Item(121, 10).Value = 1025.2174
$ws.Cells.Item(121, 12).Value = 3075.6522
$ws.Cells.Item(121, 14).Value = -6569.6522

$ws.Cells.Item(132, 8).Value = 1799.6769
$ws.Cells.Item(132, 9).Value = 1159.1017
$ws.Cells.Item(132, 10).Value = 8098.6665
$ws.Cells.Item(132, 11).Value = 3477.3051
$ws.Cells.Item(132, 12).Value = 24295.9995
$ws.Cells.Item(132, 13).Value = -947.3050999999996
$ws.Cells.Item(132, 14).Value = -29355.9995

$ws.Cells.Item(134, 8).Value = 47693.7
$ws.Cells.Item(134, 10).Value = 47693.7
$ws.Cells.Item(134, 12).Value = 47693.7
$ws.Cells.Item(134, 14).Value = -57833.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1167.25
$ws.Cells.Item(61, 9).Value = 1167.25
$ws.Cells.Item(61, 11).Value = 1167.25
$ws.Cells.Item(61, 13).Value = -955.25

$ws.Cells.Item(122, 8).Value = 68337
$ws.Cells.Item(122, 9).Value = 200012
$ws.Cells.Item(122, 10).Value = 2499.5
$ws.Cells.Item(122, 11).Value = 600036
$ws.Cells.Item(122, 12).Value = 7498.5
$ws.Cells.Item(122, 13).Value = -597586
$ws.Cells.Item(122, 14).Value = -12398.5

$ws.Cells.Item(136, 8).Value = 1167.25
$ws.Cells.Item(136, 9).Value = 1167.25
$ws.Cells.Item(136, 11).Value = 3501.75
$ws.Cells.Item(136, 13).Value = -951.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 704.8333
$ws.Cells.Item(94, 9).Value = 629.2222
$ws.Cells.Item(94, 10).Value = 931.6667
$ws.Cells.Item(94, 11).Value = 629.2222
$ws.Cells.Item(94, 12).Value = 931.6667
$ws.Cells.Item(94, 13).Value = -178.2222
$ws.Cells.Item(94, 14).Value = -1833.6667

$ws.Cells.Item(105, 8).Value = 3021.25
$ws.Cells.Item(105, 9).Value = 1786.6
$ws.Cells.Item(105, 11).Value = 1786.6
$ws.Cells.Item(105, 13).Value = -39.59999999999991

$ws.Cells.Item(134, 8).Value = 2243.2432
$ws.Cells.Item(134, 9).Value = 2074.923
$ws.Cells.Item(134, 10).Value = 2641.0908
$ws.Cells.Item(134, 11).Value = 6224.768999999999
$ws.Cells.Item(134, 12).Value = 7923.2724
$ws.Cells.Item(134, 13).Value = -3689.768999999999
$ws.Cells.Item(134, 14).Value = -12993.2724

$ws.Cells.Item(137, 8).Value = 35552
$ws.Cells.Item(137, 10).Value = 35552
$ws.Cells.Item(137, 12).Value = 35552
$ws.Cells.Item(137, 14).Value = -45752

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 33.75
$ws.Cells.Item(7, 9).Value = 30.7
$ws.Cells.Item(7, 10).Value = 49
$ws.Cells.Item(7, 11).Value = 30.7
$ws.Cells.Item(7, 12).Value = 49
$ws.Cells.Item(7, 13).Value = 82.3
$ws.Cells.Item(7, 14).Value = -275

$ws.Cells.Item(31, 8).Value = 3111.7192
$ws.Cells.Item(31, 9).Value = 1837.25
$ws.Cells.Item(31, 10).Value = 4038.606
$ws.Cells.Item(31, 11).Value = 1837.25
$ws.Cells.Item(31, 12).Value = 4038.606
$ws.Cells.Item(31, 13).Value = -1542.25
$ws.Cells.Item(31, 14).Value = -4628.606

$ws.Cells.Item(34, 8).Value = 3111.7192
$ws.Cells.Item(34, 9).Value = 1837.25
$ws.Cells.Item(34, 10).Value = 4038.606
$ws.Cells.Item(34, 11).Value = 1837.25
$ws.Cells.Item(34, 12).Value = 4038.606
$ws.Cells.Item(34, 13).Value = -1635.25
$ws.Cells.Item(34, 14).Value = -4442.606

$ws.Cells.Item(58, 8).Value = 1337.258
$ws.Cells.Item(58, 9).Value = 1316.9048
$ws.Cells.Item(58, 10).Value = 1380
$ws.Cells.Item(58, 11).Value = 1316.9048
$ws.Cells.Item(58, 12).Value = 1380
$ws.Cells.Item(58, 13).Value = -1113.9048
$ws.Cells.Item(58, 14).Value = -1786

$ws.Cells.Item(86, 8).Value = 166670340
$ws.Cells.Item(86, 9).Value = 333336500
$ws.Cells.Item(86, 10).Value = 4166
$ws.Cells.Item(86, 11).Value = 333336500
$ws.Cells.Item(86, 12).Value = 4166
$ws.Cells.Item(86, 13).Value = -333335377
$ws.Cells.Item(86, 14).Value = -6412

$ws.Cells.Item(89, 8).Value = 166670340
$ws.Cells.Item(89, 9).Value = 333336500
$ws.Cells.Item(89, 10).Value = 4166
$ws.Cells.Item(89, 11).Value = 1666682500
$ws.Cells.Item(89, 12).Value = 20830
$ws.Cells.Item(89, 13).Value = -1666676884
$ws.Cells.Item(89, 14).Value = -32062

$ws.Cells.Item(122, 8).Value = 1744.1333
$ws.Cells.Item(122, 9).Value = 1089.3846
$ws.Cells.Item(122, 10).Value = 6000
$ws.Cells.Item(122, 11).Value = 3268.1538
$ws.Cells.Item(122, 12).Value = 18000
$ws.Cells.Item(122, 13).Value = -818.1538
$ws.Cells.Item(122, 14).Value = -22900

$ws.Cells.Item(136, 8).Value = 1337.258
$ws.Cells.Item(136, 9).Value = 1316.9048
$ws.Cells.Item(136, 10).Value = 1380
$ws.Cells.Item(136, 11).Value = 3950.7144
$ws.Cells.Item(136, 12).Value = 4140
$ws.Cells.Item(136, 13).Value = -1400.7144
$ws.Cells.Item(136, 14).Value = -9240

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 3354.5454
$ws.Cells.Item(56, 9).Value = 3354.5454
$ws.Cells.Item(56, 11).Value = 3354.5454
$ws.Cells.Item(56, 13).Value = -2824.5454

$ws.Cells.Item(60, 8).Value = 1188.8889
$ws.Cells.Item(60, 10).Value = 5000
$ws.Cells.Item(60, 12).Value = 15000
$ws.Cells.Item(60, 14).Value = -15502

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 31.933332
$ws.Cells.Item(2, 9).Value = 25.75
$ws.Cells.Item(2, 10).Value = 39
$ws.Cells.Item(2, 11).Value = 25.75
$ws.Cells.Item(2, 12).Value = 39
$ws.Cells.Item(2, 13).Value = 87.25
$ws.Cells.Item(2, 14).Value = -265

$ws.Cells.Item(80, 8).Value = 2973.2666
$ws.Cells.Item(80, 9).Value = 2729.9
$ws.Cells.Item(80, 10).Value = 3460
$ws.Cells.Item(80, 11).Value = 2729.9
$ws.Cells.Item(80, 12).Value = 3460
$ws.Cells.Item(80, 13).Value = -1731.9
$ws.Cells.Item(80, 14).Value = -5456

$ws.Cells.Item(83, 8).Value = 2973.2666
$ws.Cells.Item(83, 9).Value = 2729.9
$ws.Cells.Item(83, 10).Value = 3460
$ws.Cells.Item(83, 11).Value = 13649.5
$ws.Cells.Item(83, 12).Value = 17300
$ws.Cells.Item(83, 13).Value = -8657.5
$ws.Cells.Item(83, 14).Value = -27284

$ws.Cells.Item(126, 8).Value = 10173
$ws.Cells.Item(126, 9).Value = 2479.9167
$ws.Cells.Item(126, 10).Value = 14369.228
$ws.Cells.Item(126, 11).Value = 7439.750100000001
$ws.Cells.Item(126, 12).Value = 43107.68399999999
$ws.Cells.Item(126, 13).Value = -4969.750100000001
$ws.Cells.Item(126, 14).Value = -48047.68399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 1728.9231
$ws.Cells.Item(68, 9).Value = 1706.3334
$ws.Cells.Item(68, 10).Value = 2000
$ws.Cells.Item(68, 11).Value = 1706.3334
$ws.Cells.Item(68, 12).Value = 2000
$ws.Cells.Item(68, 13).Value = -957.3334
$ws.Cells.Item(68, 14).Value = -3498

$ws.Cells.Item(71, 8).Value = 1728.9231
$ws.Cells.Item(71, 9).Value = 1706.3334
$ws.Cells.Item(71, 10).Value = 2000
$ws.Cells.Item(71, 11).Value = 8531.666999999999
$ws.Cells.Item(71, 12).Value = 10000
$ws.Cells.Item(71, 13).Value = -4787.666999999999
$ws.Cells.Item(71, 14).Value = -17488

$ws.Cells.Item(132, 8).Value = 5966.905
$ws.Cells.Item(132, 9).Value = 6221.2
$ws.Cells.Item(132, 10).Value = 5331.1665
$ws.Cells.Item(132, 11).Value = 18663.6
$ws.Cells.Item(132, 12).Value = 15993.4995
$ws.Cells.Item(132, 13).Value = -16133.6
$ws.Cells.Item(132, 14).Value = -21053.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 971.67444
$ws.Cells.Item(136, 9).Value = 1007.675
$ws.Cells.Item(136, 10).Value = 491.66666
$ws.Cells.Item(136, 11).Value = 3023.025
$ws.Cells.Item(136, 12).Value = 1474.99998
$ws.Cells.Item(136, 13).Value = -473.0249999999996
$ws.Cells.Item(136, 14).Value = -6574.999980000001
